$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "% severity levels 1-2"
$ws.Range("E1").Value = "# severity levels 1-2"
$ws.Range("F1").Value = "% severity level 3"
$ws.Range("G1").Value = "# severity level 3"
$ws.Range("H1").Value = "% severity level 4"
$ws.Range("I1").Value = "# severity level 4"
$ws.Range("J1").Value = "% severity level 5"
$ws.Range("K1").Value = "# severity level 5"
$ws.Range("L1").Value = "% Tot PiN (severity levels 3-5)"
$ws.Range("M1").Value = "# Tot PiN (severity levels 3-5)"
